$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Natmi LR-pairs table update ("following Dr Hou advice"): a new target
# cluster "ECs" is added to the Wnt5a -> Fzd4 signalling summary, and the
# per-row stats are recomputed across the full Sending x Target cluster
# cross-join (FAPs/sCs senders x FAPs/sCs/ECs targets).
# ---------------------------------------------------------------------------

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 22.41709
$ws.Range("N2").Value = 67.25127
$ws.Range("O2").Value = 0.3988455747018376
$ws.Range("P2").Value = 0.3988455747018376
$ws.Range("Q2").Value = 205.2723815016734
$ws.Range("R2").Value = 1847.45143351506
$ws.Range("S2").Value = 0.3866686054309541
$ws.Range("T2").Value = 0.3866686054309541

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.78189033333333
$ws.Range("N3").Value = 50.345671
$ws.Range("O3").Value = 0.2985839238983091
$ws.Range("P3").Value = 0.2985839238983091
$ws.Range("Q3").Value = 153.6710873187931
$ws.Range("R3").Value = 1383.039785869138
$ws.Range("S3").Value = 0.289467996590334
$ws.Range("T3").Value = 0.289467996590334

# Row 4: FAPs -> sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.00595566666667
$ws.Range("N4").Value = 51.017867
$ws.Range("O4").Value = 0.3025705013998533
$ws.Range("P4").Value = 0.3025705013998533
$ws.Range("Q4").Value = 155.7228444641362
$ws.Range("R4").Value = 1401.505600177226
$ws.Range("S4").Value = 0.2933328617430108
$ws.Range("T4").Value = 0.2933328617430108

# Row 5: sCs -> ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 22.41709
$ws.Range("N5").Value = 67.25127
$ws.Range("O5").Value = 0.3988455747018376
$ws.Range("P5").Value = 0.3988455747018376
$ws.Range("Q5").Value = 6.46443866039
$ws.Range("R5").Value = 58.17994794351001
$ws.Range("S5").Value = 0.01217696927088348
$ws.Range("T5").Value = 0.01217696927088348

# Row 6: sCs -> FAPs (new row)
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.78189033333333
$ws.Range("N6").Value = 50.345671
$ws.Range("O6").Value = 0.2985839238983091
$ws.Range("P6").Value = 0.2985839238983091
$ws.Range("Q6").Value = 4.839410497313667
$ws.Range("R6").Value = 43.554694475823
$ws.Range("S6").Value = 0.009115927307975141
$ws.Range("T6").Value = 0.009115927307975143

# Row 7: sCs -> sCs (new row)
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.00595566666667
$ws.Range("N7").Value = 51.017867
$ws.Range("O7").Value = 0.3025705013998533
$ws.Range("P7").Value = 0.3025705013998533
$ws.Range("Q7").Value = 4.904024441552333
$ws.Range("R7").Value = 44.136219973971
$ws.Range("S7").Value = 0.009237639656842466
$ws.Range("T7").Value = 0.009237639656842468
